$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Add new column S (year 2023) data ---
# Year header row
$ws.Range("S3").Value = 2023

# Data rows 4-14, column S
$ws.Range("S4").Value = 1132.8
$ws.Range("S5").Value = 182.9
$ws.Range("S6").Value = 6970
$ws.Range("S7").Value = 4164
$ws.Range("S8").Value = 733.5
$ws.Range("S9").Value = 36
$ws.Range("S10").Value = 37.1
$ws.Range("S11").Value = 147
$ws.Range("S12").Value = 7.9
$ws.Range("S13").Value = 999.1
$ws.Range("S14").Value = 965

# Copy formatting from column R to column S, row by row (skip blank row 2 so we
# don't introduce a stray cell there).
$ws.Range("R1").Copy()
$ws.Range("S1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("R3:R14").Copy()
$ws.Range("S3:S14").PasteSpecial(-4122)  # xlPasteFormats

# --- Update merged cell range A1:R1 -> A1:S1 ---
$ws.Range("A1:R1").UnMerge()
$ws.Range("A1:S1").Merge()

# --- Update sheet view: topLeftCell + selection ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("A1:S1").Select()

# --- Update workbook window size ---
$excel.ActiveWindow.Width = 22215
$excel.ActiveWindow.Height = 7305
